# "updated classes and docs"
#
# Adds a new ScheduleName / AstroSchedule row to the Configuration sheet,
# re-fits column B to the new (wider) content, and leaves Configuration as
# the active sheet/tab with B8 selected (Entries consequently loses the
# "active tab" status it previously had).

$wb = $excel.ActiveWorkbook

$config = $wb.Worksheets.Item("Configuration")

# New configuration entry appended after the existing
# EBOVersion / ReferenceYear / Comments rows.
$config.Range("A7").Value = "ScheduleName"
$config.Range("B7").Value = "AstroSchedule"

# "AstroSchedule" is wider than the previous contents of column B, so
# re-fit the column to the new data.
$config.Columns.Item(2).AutoFit() | Out-Null

# Switch the active sheet/tab to Configuration and leave the selection on B8.
$config.Activate()
$config.Range("B8").Select() | Out-Null
